$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM row 11: LED1 / 0805 / C2297 (mirrors the style of the existing rows)
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Rows("11").RowHeight = $ws.Rows("10").RowHeight

$ws.Range("B11").Value = "LED1"
$ws.Range("C11").Value = "'0805"
$ws.Range("D11").Value = "C2297"

$ws.Application.CutCopyMode = 0

$ws.Range("D12").Select() | Out-Null
